$d = $word.ActiveDocument

$ids = @("p065v_1", "p065v_2", "p065v_3", "p065v_4", "p065v_5", "p065v_6", "p065v_7")

foreach ($id in $ids) {
    $old = "<id>" + $id + "</id>"
    $new = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
